$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.359.71'
$ws.Range("E2").Value = '  -2.02%  '

$ws.Range("D3").Value = '1.796.24'
$ws.Range("E3").Value = '  -1.79%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.006'
$ws.Range("E5").Value = '  -0.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '307.49'
$ws.Range("E6").Value = '  -1.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4509'
$ws.Range("E7").Value = '  -1.66%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3595'
$ws.Range("E8").Value = '  -2.65%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.01'
$ws.Range("E9").Value = '  +0.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07080'
$ws.Range("E10").Value = '  -1.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8851'
$ws.Range("E11").Value = '  +0.92%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07806'
$ws.Range("E12").Value = '  -0.39%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.42'
$ws.Range("E13").Value = '  -1.01%  '

$ws.Range("D14").Value = '1.774.21'
$ws.Range("E14").Value = '  -3.06%  '

$ws.Range("E15").Value = '  -0.89%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.351'
$ws.Range("E16").Value = '  -0.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '84.92'
$ws.Range("E17").Value = '  -2.64%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.008'
$ws.Range("E18").Value = '  -0.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008515'
$ws.Range("E19").Value = '  -2.24%  '

$ws.Range("E20").Value = '  -0.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.27'
$ws.Range("E21").Value = '  -1.40%  '

$ws.Range("D22").Value = '26.370.28'
$ws.Range("E22").Value = '  -2.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.992'
$ws.Range("E23").Value = '  -0.13%  '

$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.54'
$ws.Range("E24").Value = '  +0.96%  '

$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D25").Value = '2.006.63'
$ws.Range("E25").Value = '  -2.52%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.977'
$ws.Range("E26").Value = '  +0.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.32'
$ws.Range("E27").Value = '  +0.90%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.85'
$ws.Range("E28").Value = '  -1.75%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.034'
$ws.Range("E29").Value = '  +3.30%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '112.12'
$ws.Range("E30").Value = '  -1.56%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.876'
$ws.Range("E31").Value = '  -0.94%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08684'
$ws.Range("E32").Value = '  -1.32%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.084'
$ws.Range("E33").Value = '  +1.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.754'
$ws.Range("E34").Value = '  +7.35%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.450'
$ws.Range("E35").Value = '  -0.74%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7245'
$ws.Range("E36").Value = '  -3.92%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.107'
$ws.Range("E37").Value = '  -2.30%  '

$ws.Range("E39").Value = '  -1.41%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01928'
$ws.Range("E40").Value = '  -0.31%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.05093'
$ws.Range("E41").Value = '  -0.80%  '

$ws.Range("E42").Value = '  -1.33%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.889'
$ws.Range("E43").Value = '  -0.81%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5073'
$ws.Range("E44").Value = '  +1.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1513'
$ws.Range("E45").Value = '  -5.13%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.011'
$ws.Range("E46").Value = '  -3.49%  '

$ws.Range("E47").Value = '  -0.12%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4626'
$ws.Range("E48").Value = '  -1.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.32'
$ws.Range("E49").Value = '  -0.77%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.860'
$ws.Range("E50").Value = '  -3.04%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.579'
